$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = -0.4340006244244006
$ws.Range("J2").Value = 0.2311415746764876
$ws.Range("K2").Value = -0.2338051555780229
$ws.Range("L2").Value = 2.70624499875677
